$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OCS Input")

# Update Well Name (B5): GK-W527B -> GK-P527B
$ws.Range("B5").Value = "GK-P527B"

# Update Charging Mechanism for row 12 (Item 2)
$ws.Range("L12").Value = "0.5 unit/day on end phase 15"

# Update Charging Mechanism for row 11 (Item 1)
$ws.Range("L11").Value = "1.2 unit/day on 2023/01/03"

# Update selection to reflect last active cell L12
$ws.Range("L12").Select()
